# Refresh "cryptos" price/volume data (GitHub Actions scheduled run).
# Column D ("Price") holds number-like text (dot-grouped, e.g. "28.467.21")
# and column E ("Volume(1h)") holds padded percentage text; both must stay
# text cells, so values that Excel would otherwise auto-convert to a number
# are entered with a leading apostrophe (same as typing them in the UI).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.467.21"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "1.863.33"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'324.74"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.4554"
$ws.Range("E7").Value = "  -2.21%  "
$ws.Range("D8").Value = "'0.3829"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "'0.07812"
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").Value = "'0.9845"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D11").Value = "'21.46"
$ws.Range("E11").Value = "  -4.06%  "
$ws.Range("D12").Value = "1.872.85"
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("D13").Value = "'6.897"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "'5.633"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "'0.06914"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D17").Value = "'86.49"
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "'16.68"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "28.472.32"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").Value = "'5.249"
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("D23").Value = "'10.88"
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("D24").Value = "'2.088"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").Value = "2.082.54"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").Value = "'153.27"
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").Value = "'5.657"
$ws.Range("E28").Value = "  -2.25%  "
$ws.Range("D29").Value = "'117.19"
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("D30").Value = "'1.884"
$ws.Range("E30").Value = "  -5.82%  "
$ws.Range("D31").Value = "'0.09267"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").Value = "'0.9038"
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("D33").Value = "'5.275"
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("D35").Value = "'3.287"
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("D36").Value = "'0.05664"
$ws.Range("E36").Value = "  -3.12%  "
$ws.Range("D37").Value = "'1.144"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "'0.02034"
$ws.Range("E38").Value = "  -4.07%  "
$ws.Range("D39").Value = "'7.622"
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("D40").Value = "'0.5547"
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("D41").Value = "'0.1763"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").Value = "'9.576"
$ws.Range("D43").Value = "'0.07134"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("D45").Value = "'0.5226"
$ws.Range("D46").Value = "'1.123"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("D47").Value = "'2.104"
$ws.Range("E47").Value = "  -7.03%  "
$ws.Range("D48").Value = "'1.803"
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("D49").Value = "'111.68"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("E51").Value = "  +0.12%  "
